# Atualização automática de CANDELARIA.xlsx
#
# Changes applied:
#   1. Rename sheet "Paineis DARQ" -> "PAINEIS DARQ"
#   2. Rename sheet "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO"
#   3. Delete sheet "Desarquivamentos Pendentes" entirely
#
# (The "DGC" sheet keeps its data/position; its sheetId/rId shift automatically
#  once the preceding sheet is removed.)

$wb = $excel.ActiveWorkbook

# Avoid any interactive "are you sure you want to delete this sheet" prompt.
$excel.DisplayAlerts = $false | Out-Null

# 1. Rename "Paineis DARQ" -> "PAINEIS DARQ"
$wsPaineis = $wb.Worksheets.Item("Paineis DARQ")
$wsPaineis.Name = "PAINEIS DARQ"

# 2. Rename "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO"
$wsRecolhimento = $wb.Worksheets.Item("Recolhimento x Eliminacao")
$wsRecolhimento.Name = "RECOLHIMENTO X ELIMINAÇÃO"

# 3. Delete the "Desarquivamentos Pendentes" sheet
$wsDesarquivamentos = $wb.Worksheets.Item("Desarquivamentos Pendentes")
$wsDesarquivamentos.Delete() | Out-Null

$excel.DisplayAlerts = $true | Out-Null

Write-Host "Sheets after edit:"
foreach ($s in $wb.Worksheets) {
    Write-Host (" - " + $s.Name)
}
